# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计" with the
#   quarter's fund-holdings detail (same column layout as "2021-Q4").
# - Insert a new leading row into "总计" summarising the 2022-Q1 totals,
#   pushing the existing "2021-Q4" summary row down.

$wb = $excel.ActiveWorkbook

$wsQ4 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, positioned right after "2021-Q4"
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

# NOTE: sheet handles obtained before a Worksheets.Add() become stale
# positional anchors once the tab order shifts, so re-fetch "总计" only
# after the insert is done (see below).

# Header row (same headings/style as the "2021-Q4" sheet)
$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)

$wsQ1.Cells.Item(1,2).Value = "基金代码"
$wsQ1.Cells.Item(1,3).Value = "基金名称"
$wsQ1.Cells.Item(1,4).Value = "基金规模"
$wsQ1.Cells.Item(1,5).Value = "股票总仓位"
$wsQ1.Cells.Item(1,6).Value = "仓位占比"
$wsQ1.Cells.Item(1,7).Value = "持有市值(亿元)"
$wsQ1.Cells.Item(1,8).Value = "仓位排名"

# Row-index column (A) style, copied down for all 8 data rows
$wsQ4.Range("A2").Copy()
$wsQ1.Range("A2:A9").PasteSpecial(-4122)

# Columns B,C,D,E,F,G hold text (fund code/name/size/position/%/value are
# all stored as text in this workbook, even the numeric-looking ones) -
# force text formatting so numeric-looking strings (fund codes, sizes,
# percentages) are not auto-converted to numbers.
$wsQ1.Range("B2:G9").NumberFormat = "@"

$data = @(
    @("002446","广发利鑫灵活配置混合A","12.46","74.35","2.88","0.3588",10),
    @("001471","融通新能源灵活配置混合","5.77","83.30","3.29","0.1898",9),
    @("011172","广发利鑫灵活配置混合C","1.10","74.35","2.88","0.0317",10),
    @("006522","财通新兴蓝筹混合A","0.29","90.33","4.10","0.0119",6),
    @("006890","上投摩根领先优选混合","0.36","79.50","3.22","0.0116",7),
    @("001830","融通跨界成长灵活配置混合","0.25","77.20","2.65","0.0066",8),
    @("620002","金元顺安成长动力混合","0.16","62.78","2.94","0.0047",7),
    @("006523","财通新兴蓝筹混合C","0.03","90.33","4.10","0.0012",6)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    $wsQ1.Cells.Item($row,1).Value = $i
    $wsQ1.Cells.Item($row,2).Value = $rec[0]
    $wsQ1.Cells.Item($row,3).Value = $rec[1]
    $wsQ1.Cells.Item($row,4).Value = $rec[2]
    $wsQ1.Cells.Item($row,5).Value = $rec[3]
    $wsQ1.Cells.Item($row,6).Value = $rec[4]
    $wsQ1.Cells.Item($row,7).Value = $rec[5]
    $wsQ1.Cells.Item($row,8).Value = $rec[6]
}

# ---------------------------------------------------------------------
# 2) "总计" sheet: insert a new first data row for 2022-Q1, pushing the
#    existing 2021-Q4 summary row down to row 3.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

# Column-A "index" style (bold/centered, same as used in both sheets'
# A columns) needs to be (re)applied explicitly to the new row - Insert()
# does not carry it over onto the freshly created A2 cell.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q1"
$wsTotal.Cells.Item(2,3).Value = 8
$wsTotal.Cells.Item(2,4).Value = 0.62

$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(3,2).Value = "2021-Q4"
$wsTotal.Cells.Item(3,3).Value = 2
$wsTotal.Cells.Item(3,4).Value = 0.06
